# Auto-update draw results: append the 2025-09-26 Pick 4 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 10

# Columns A (date) and C (phase code) look numeric/date-like to Excel's
# auto-detection, but the source data stores them as plain text (matching
# every prior row in the sheet) — force text format before assigning so
# they aren't reinterpreted as a date serial / number.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "2025-09-26"

$ws.Range("B$newRow").Value = "Pick 4"

$ws.Range("C$newRow").NumberFormat = "@"
$ws.Range("C$newRow").Value = "250926"

$ws.Range("D$newRow").Value = "9-6-0-3"

$ws.Range("E$newRow").Value = "2025-09-26T21:37:23.252+04:00"
